# The commit swaps the presentation's theme (ppt/theme/theme1.xml, "Integral")
# for the stock default "Office Theme" color palette that used to live in
# ppt/theme/theme2.xml (only ever referenced by the Notes Master).
#
# The 12 DrawingML theme colors (dk1, lt1, dk2, lt2, accent1-6, hlink,
# folHlink) are reachable/writable on every Slide's ThemeColorScheme
# collection (they all point at the single shared slide-master theme part).
# Re-point each of them at the "Office" theme's RGB values.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# index -> (theme color slot, target "Office" RGB)
$tcs.Item(1).RGB  = 0        # dk1      000000
$tcs.Item(2).RGB  = 16777215 # lt1      FFFFFF
$tcs.Item(3).RGB  = 6968388  # dk2      44546A
$tcs.Item(4).RGB  = 15132391 # lt2      E7E6E6
$tcs.Item(5).RGB  = 13998939 # accent1  5B9BD5
$tcs.Item(6).RGB  = 3243501  # accent2  ED7D31
$tcs.Item(7).RGB  = 10855845 # accent3  A5A5A5
$tcs.Item(8).RGB  = 49407    # accent4  FFC000
$tcs.Item(9).RGB  = 12874308 # accent5  4472C4
$tcs.Item(10).RGB = 4697456  # accent6  70AD47
$tcs.Item(11).RGB = 12673797 # hlink    0563C1
$tcs.Item(12).RGB = 7491477  # folHlink 954F72
